$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# ---- 1. Geometry (off/ext) updates on existing shapes ----
$sh = Get-ShapeById $s 118
$sh.Width = 595.82165527343750000000

$sh = Get-ShapeById $s 76
$sh.Left = 612.16333007812500000000

$sh = Get-ShapeById $s 79
$sh.Left = 562.99645996093750000000
$sh.Top = 213.15985107421875000000
$sh.Width = 49.16685104370117187500
$sh.Height = 25.78378105163574218750

$sh = Get-ShapeById $s 80
$sh.Left = 612.16333007812500000000
$sh.Top = 227.69213867187500000000

$sh = Get-ShapeById $s 81
$sh.Top = 238.94346618652343750000
$sh.Width = 39.09283828735351562500
$sh.Height = 0.02433070912957191467

$sh = Get-ShapeById $s 83
$sh.Left = 612.16333007812500000000

$sh = Get-ShapeById $s 84
$sh.Width = 39.09283828735351562500

$sh = Get-ShapeById $s 85
$sh.Left = 612.16333007812500000000

$sh = Get-ShapeById $s 86
$sh.Width = 39.09283828735351562500

$sh = Get-ShapeById $s 52
$sh.Left = 612.16333007812500000000

$sh = Get-ShapeById $s 53
$sh.Width = 39.09283828735351562500

$sh = Get-ShapeById $s 56
$sh.Left = 612.00000000000000000000

$sh = Get-ShapeById $s 57
$sh.Width = 38.92953109741210937500

$sh = Get-ShapeById $s 59
$sh.Left = 612.00000000000000000000

$sh = Get-ShapeById $s 67
$sh.Width = 38.92953109741210937500

# ---- 2. bentConnector3 adj1 guide updates ----
$sh = Get-ShapeById $s 79
$sh.Adjustments.Item(1) = 0.50073003768920898438

$sh = Get-ShapeById $s 84
$sh.Adjustments.Item(1) = 0.36750000715255737305

$sh = Get-ShapeById $s 86
$sh.Adjustments.Item(1) = 0.36750000715255737305

$sh = Get-ShapeById $s 53
$sh.Adjustments.Item(1) = 0.36750000715255737305

$sh = Get-ShapeById $s 57
$sh.Adjustments.Item(1) = 0.36694002151489257812

$sh = Get-ShapeById $s 67
$sh.Adjustments.Item(1) = 0.36694002151489257812

# ---- 3. Delete the old '*' multiplicity textboxes ----
$sh = Get-ShapeById $s 72
$sh.Delete()
$sh = Get-ShapeById $s 73
$sh.Delete()
$sh = Get-ShapeById $s 74
$sh.Delete()
$sh = Get-ShapeById $s 77
$sh.Delete()

# ---- 4. Update the remaining old textbox (id 82) in place ----
$sh = Get-ShapeById $s 82
$sh.Left = 587.21588134765625000000
$sh.Top = 226.81410217285156250000
$sh.Width = 19.33047294616699218750
$sh.Height = 14.06952762603759765625
$tr = $sh.TextFrame.TextRange
$tr.Text = '0..1'
$tr.Font.Size = 10.0

# ---- 5. Add the new multiplicity textboxes ----
$tb = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tb.Left = 587.21588134765625000000
$tb.Top = 252.34851074218750000000
$tb.Width = 19.33047294616699218750
$tb.Height = 14.06952762603759765625
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 0
$tb.TextFrame.MarginLeft = 0
$tb.TextFrame.MarginTop = 0
$tb.TextFrame.MarginRight = 0
$tb.TextFrame.MarginBottom = 0
$tb.Fill.Visible = 0
$tr = $tb.TextFrame.TextRange
$tr.Text = '0..1'
$tr.Font.Size = 10.0
$tr.Font.Color.RGB = 10498160
$tr.ParagraphFormat.Alignment = 2
$tb.Name = 'TextBox 70'

$tb = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tb.Left = 587.89318847656250000000
$tb.Top = 276.00000000000000000000
$tb.Width = 19.33047294616699218750
$tb.Height = 14.06952762603759765625
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 0
$tb.TextFrame.MarginLeft = 0
$tb.TextFrame.MarginTop = 0
$tb.TextFrame.MarginRight = 0
$tb.TextFrame.MarginBottom = 0
$tb.Fill.Visible = 0
$tr = $tb.TextFrame.TextRange
$tr.Text = '0..1'
$tr.Font.Size = 10.0
$tr.Font.Color.RGB = 10498160
$tr.ParagraphFormat.Alignment = 2
$tb.Name = 'TextBox 86'

$tb = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tb.Left = 587.89318847656250000000
$tb.Top = 304.66806030273437500000
$tb.Width = 19.33047294616699218750
$tb.Height = 14.06952762603759765625
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 0
$tb.TextFrame.MarginLeft = 0
$tb.TextFrame.MarginTop = 0
$tb.TextFrame.MarginRight = 0
$tb.TextFrame.MarginBottom = 0
$tb.Fill.Visible = 0
$tr = $tb.TextFrame.TextRange
$tr.Text = '0..1'
$tr.Font.Size = 10.0
$tr.Font.Color.RGB = 10498160
$tr.ParagraphFormat.Alignment = 2
$tb.Name = 'TextBox 87'

$tb = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tb.Left = 587.21575927734375000000
$tb.Top = 330.20245361328125000000
$tb.Width = 19.33047294616699218750
$tb.Height = 14.06952762603759765625
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 0
$tb.TextFrame.MarginLeft = 0
$tb.TextFrame.MarginTop = 0
$tb.TextFrame.MarginRight = 0
$tb.TextFrame.MarginBottom = 0
$tb.Fill.Visible = 0
$tr = $tb.TextFrame.TextRange
$tr.Text = '0..1'
$tr.Font.Size = 10.0
$tr.Font.Color.RGB = 10498160
$tr.ParagraphFormat.Alignment = 2
$tb.Name = 'TextBox 88'

$tb = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tb.Left = 586.43884277343750000000
$tb.Top = 201.93048095703125000000
$tb.Width = 19.33047294616699218750
$tb.Height = 14.06952762603759765625
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 0
$tb.TextFrame.MarginLeft = 0
$tb.TextFrame.MarginTop = 0
$tb.TextFrame.MarginRight = 0
$tb.TextFrame.MarginBottom = 0
$tb.Fill.Visible = 0
$tr = $tb.TextFrame.TextRange
$tr.Text = '1'
$tr.Font.Size = 10.0
$tr.Font.Color.RGB = 10498160
$tr.ParagraphFormat.Alignment = 2
$tb.Name = 'TextBox 71'

